# Add a third "tags" column (comma separated keyword list) to each of the
# five sheets, in column C, row 2 (next to the first data row under the
# header). Shared strings must be appended in the order: Sheet3 (hoodie),
# Sheet1 (polo), Sheet2 (fullsleeve), Sheet5 (jersey), Sheet4 (coat) so the
# new shared-string table lines up with the target workbook.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)
$ws5 = $wb.Worksheets.Item(5)

# --- write the new tag values (this also drives shared-string order) ---
$ws3.Range("C2").Value = "nike,blue,tshirt,pro,white,hoodie,star,grey,nik,green,yo,yellow,adidas,red"
$ws1.Range("C2").Value = "nike,blue,tshirt,pro,white,polo,star,grey,nik,green,yo,yellow,adidas,red"
$ws2.Range("C2").Value = "nike,blue,tshirt,pro,white,fullsleeve,star,grey,nik,green,yo,yellow,adidas,red"
$ws5.Range("C2").Value = "nike,blue,tshirt,pro,white,jersey,star,grey,nik,green,yo,yellow,adidas,red"
$ws4.Range("C2").Value = "nike,blue,tshirt,pro,white,coat,star,grey,nik,green,yo,yellow,adidas,red"

# --- widen column C on each sheet to fit the new tag text ---
$ws1.Columns.Item(3).ColumnWidth = 77.3333333333333
$ws2.Columns.Item(3).ColumnWidth = 91.1666666666667
$ws3.Columns.Item(3).ColumnWidth = 74.3333333333333
$ws4.Columns.Item(3).ColumnWidth = 72.1666666666667
$ws5.Columns.Item(3).ColumnWidth = 72.8333333333333

# --- move the selection to the new cell on each sheet ---
$ws1.Range("C2").Select()
$ws2.Range("C2").Select()
$ws4.Range("C2").Select()

$ws3.Range("C4").Select()

# Re-activate Sheet5 last (it was, and remains, the selected tab) and park
# the selection on its new C2 cell.
$ws5.Activate()
$ws5.Range("C2").Select()
